# Automatische test-sync: 2025-08-26 21:12:50
#
# Adds the newest "Retour status" mail-log entry (2025-08-26 21:11:53) to the
# "Logs" sheet as row 10, extends the conditional formatting ranges that used
# to stop at row 9 so they now cover row 10 as well, and refreshes the
# "Dashboard" summary sheet so that it reflects the updated category counts
# (Retour / Terugbetaling now has 5 hits, Klantenservice / Opvolging keeps 4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Logs sheet: append the new row of data
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(10, 1).Value  = "Retour status"
$logs.Cells.Item(10, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(10, 4).Value  = "Retour / Terugbetaling"
$logs.Cells.Item(10, 6).Value  = "2025-08-26 21:11:53"
$logs.Cells.Item(10, 7).Value  = "Ja"
$logs.Cells.Item(10, 8).Value  = "Nee"
$logs.Cells.Item(10, 9).Value  = "Nee"
$logs.Cells.Item(10, 10).Value = "Nee"

# ---------------------------------------------------------------------------
# 2) Logs sheet: extend the conditional formatting ranges from row 9 to row 10
#    (each block's rules all move together when the first rule's applies-to
#    range is modified).
# ---------------------------------------------------------------------------
$logs.Range("D2:D9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D10"))
$logs.Range("G2:G9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G10"))
$logs.Range("H2:H9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H10"))
$logs.Range("I2:I9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I10"))
$logs.Range("J2:J9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J10"))

# ---------------------------------------------------------------------------
# 3) Dashboard sheet: resort/refresh the category summary so row 2 is the
#    category with the most hits ("Retour / Terugbetaling" = 5) and row 3 is
#    "Klantenservice / Opvolging" = 4.
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Retour / Terugbetaling"
$dashboard.Range("B2").Value = 5
$dashboard.Range("A3").Value = "Klantenservice / Opvolging"
$dashboard.Range("B3").Value = 4
